$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to Text
# format before assignment, then the format is reset to Normal style so the
# cell keeps matching its original (unstyled) appearance while the stored
# value remains a text string, matching the source data.
$numericRiskCells = @(
    "D5", "D6", "D7", "D9", "D10", "D11", "D12", "D13", "D15", "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D28", "D31", "D32", "D33", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48"
)
foreach ($addr in $numericRiskCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = '42.418.08'
$ws.Range("E2").Value = '  -1.65%  '
$ws.Range("D3").Value = '2.186.53'
$ws.Range("E3").Value = '  -2.48%  '
$ws.Range("E4").Value = '  -0.41%  '
$ws.Range("D5").Value = '252.14'
$ws.Range("E5").Value = '  +2.53%  '
$ws.Range("D6").Value = '0.614'
$ws.Range("E6").Value = '  -0.88%  '
$ws.Range("D7").Value = '75.02'
$ws.Range("E7").Value = '  -1.27%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").Value = '0.582'
$ws.Range("E9").Value = '  -5.53%  '
$ws.Range("D10").Value = '40.32'
$ws.Range("E10").Value = '  -2.28%  '
$ws.Range("D11").Value = '0.0912'
$ws.Range("E11").Value = '  -2.53%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").Value = '6.81'
$ws.Range("E12").Value = '  -2.55%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '0.101'
$ws.Range("E13").Value = '  -0.12%  '
$ws.Range("D14").Value = '2.512.36'
$ws.Range("E14").Value = '  -1.81%  '
$ws.Range("D15").Value = '14.19'
$ws.Range("E15").Value = '  -3.72%  '
$ws.Range("D16").Value = '2.193.39'
$ws.Range("E16").Value = '  -1.73%  '
$ws.Range("E17").Value = '  -5.65%  '
$ws.Range("D18").Value = '42.358.82'
$ws.Range("E18").Value = '  -1.59%  '
$ws.Range("E19").Value = '  -3.57%  '
$ws.Range("D20").Value = '70.71'
$ws.Range("E20").Value = '  -0.61%  '
$ws.Range("D21").Value = '5.88'
$ws.Range("E21").Value = '  -2.11%  '
$ws.Range("D22").Value = '225.40'
$ws.Range("E22").Value = '  -2.34%  '
$ws.Range("B23").Value = 'ImmutableX'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D23").Value = '2.13'
$ws.Range("E23").Value = '  -3.77%  '
$ws.Range("B24").Value = 'InternetComputer(DFINITY)'
$ws.Range("C24").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D24").Value = '9.40'
$ws.Range("E24").Value = '  -10.25%  '
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("D26").Value = '10.46'
$ws.Range("E26").Value = '  -4.80%  '
$ws.Range("D27").Value = '3.36'
$ws.Range("E27").Value = '  -0.30%  '
$ws.Range("D28").Value = '38.45'
$ws.Range("E28").Value = '  +1.99%  '
$ws.Range("E29").Value = '  -0.84%  '
$ws.Range("E30").Value = '  -4.49%  '
$ws.Range("D31").Value = '172.84'
$ws.Range("E31").Value = '  -0.73%  '
$ws.Range("D32").Value = '20.06'
$ws.Range("E32").Value = '  -1.56%  '
$ws.Range("D33").Value = '0.0824'
$ws.Range("E33").Value = '  +3.76%  '
$ws.Range("E34").Value = '  -5.34%  '
$ws.Range("E35").Value = '  -1.79%  '
$ws.Range("D36").Value = '0.108'
$ws.Range("E36").Value = '  -3.22%  '
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").Value = '4.23'
$ws.Range("E37").Value = '  -3.17%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.0337'
$ws.Range("E38").Value = '  +1.19%  '
$ws.Range("D39").Value = '11.95'
$ws.Range("E39").Value = '  -9.22%  '
$ws.Range("D40").Value = '2.07'
$ws.Range("E40").Value = '  -3.69%  '
$ws.Range("D41").Value = '2.59'
$ws.Range("E41").Value = '  +11.80%  '
$ws.Range("B42").Value = 'MultiversX'
$ws.Range("C42").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D42").Value = '58.90'
$ws.Range("E42").Value = '  -2.44%  '
$ws.Range("B43").Value = 'THORChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D43").Value = '5.16'
$ws.Range("E43").Value = '  -7.72%  '
$ws.Range("E44").Value = '  -3.56%  '
$ws.Range("D45").Value = '101.55'
$ws.Range("E45").Value = '  -3.76%  '
$ws.Range("D46").Value = '0.0972'
$ws.Range("E46").Value = '  -2.34%  '
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").Value = '8.21'
$ws.Range("E47").Value = '  -4.77%  '
$ws.Range("B48").Value = 'WOONetwork'
$ws.Range("C48").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D48").Value = '0.459'
$ws.Range("E48").Value = '  +3.74%  '
$ws.Range("E49").Value = '  -1.71%  '
$ws.Range("E50").Value = '  -2.23%  '
$ws.Range("E51").Value = '  -0.89%  '

# Restore default (Normal) style on the forced-text cells so no residual
# number-format styling is left behind on them
foreach ($addr in $numericRiskCells) {
    $ws.Range($addr).Style = "Normal"
}
